# Applies the resume text edits described by the commit diff.
$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output "WARNING: find failed for: $findText"
    }
    return $ok
}

# 1. "...using HTML Canvas which was 7x faster..." -> add comma before "which"
Replace-Text "using HTML Canvas which was 7x faster than the previous generator." "using HTML Canvas, which was 7x faster than the previous generator."

# 2. "Worked on improving ... tools like for better quality." -> "Improved ... tools for better quality."
Replace-Text "Worked on improving user flow and UI design that follows usability and accessibility best practices while utilizing auditing tools like for better quality." "Improved user flow and UI design that follows usability and accessibility best practices while utilizing auditing tools for better quality."

# 3. "Used Git as a version control system. Submitted, reviewed, and merged..." -> "...system and submitted, reviewed, and merged..."
Replace-Text "Used Git as a version control system. Submitted, reviewed, and merged pull requests on GitHub." "Used Git as a version control system and submitted, reviewed, and merged pull requests on GitHub."

# 4. "Developed an entire full-stack web application alone that involves planning, designing, implementing, and maintaining a web application (offering Website Hosting and NFT Utilities)." -> "Developed a full-stack web application that offers Website Hosting and NFT Utilities."
Replace-Text "Developed an entire full-stack web application alone that involves planning, designing, implementing, and maintaining a web application (offering Website Hosting and NFT Utilities)." "Developed a full-stack web application that offers Website Hosting and NFT Utilities."

# 5. "Built an entire website hosting platform using Next.js..." -> "Built a website hosting platform using Next.js..."
Replace-Text "Built an entire website hosting platform using Next.js" "Built a website hosting platform using Next.js"

# 6. "Built an entire full-stack web application starter pack using Next.js, TypeScript, " -> "Built a full-stack web application starter pack using Next.js, TypeScript, "
#    (trailing space trimmed from both find/replace text -- this Find
#    implementation does not match a trailing space at the end of the needle)
Replace-Text "Built an entire full-stack web application starter pack using Next.js, TypeScript," "Built a full-stack web application starter pack using Next.js, TypeScript,"

# 7. "Integrated Prisma with MongoDB for strong type-safety which resulted in better" -> "...type-safety, resulting in better"
Replace-Text "Integrated Prisma with MongoDB for strong type-safety which resulted in better" "Integrated Prisma with MongoDB for strong type-safety, resulting in better"

# 8. "Set up code formatter and git hooks using Prettier and Husky which enhanced code quality by 80% resulting in better performance." -> "...Husky, which enhanced code quality by 80% and improved performance."
Replace-Text "Set up code formatter and git hooks using Prettier and Husky which enhanced code quality by 80% resulting in better performance." "Set up code formatter and git hooks using Prettier and Husky, which enhanced code quality by 80% and improved performance."

# 9. Remove the extra empty "List Paragraph"-styled paragraph that trails the
#    "...database." bullet (there were two blank paragraphs back-to-back; one
#    of them -- the ListParagraph-styled one -- was deleted).
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "database\.\s*$") {
        $dup = $d.Paragraphs.Item($i + 1)
        if ($dup.Style.NameLocal -eq "List Paragraph" -and $dup.Range.Text.Trim().Length -eq 0) {
            $dup.Range.Delete()
        }
        break
    }
}
